$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = "29.476.71"
$ws.Range("E2").Value = "  +1.84%  "

# Row 3
$ws.Range("D3").Value = "1.984.55"
$ws.Range("E3").Value = "  +5.00%  "

# Row 4
$ws.Range("D4").Value = "'0.9958"
$ws.Range("E4").Value = "  -0.66%  "

# Row 5
$ws.Range("D5").Value = "'327.56"
$ws.Range("E5").Value = "  +0.47%  "

# Row 6
$ws.Range("D6").Value = "'0.9980"
$ws.Range("E6").Value = "  -0.35%  "

# Row 7
$ws.Range("E7").Value = "  +1.56%  "

# Row 8
$ws.Range("D8").Value = "'0.3952"
$ws.Range("E8").Value = "  +1.23%  "

# Row 9
$ws.Range("B9").Value = "OKB"
$ws.Range("C9").Value = "https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"
$ws.Range("D9").Value = "'46.22"
$ws.Range("E9").Value = "  -1.18%  "

# Row 10
$ws.Range("B10").Value = "Dogecoin"
$ws.Range("C10").Value = "https://coinranking.com/coin/a91GCGd_u96cF+dogecoin-doge"
$ws.Range("D10").Value = "'0.07913"
$ws.Range("E10").Value = "  +0.79%  "

# Row 11
$ws.Range("B11").Value = "Polygon"
$ws.Range("C11").Value = "https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic"
$ws.Range("D11").Value = "'1.002"
$ws.Range("E11").Value = "  +1.35%  "

# Row 12
$ws.Range("B12").Value = "Solana"
$ws.Range("C12").Value = "https://coinranking.com/coin/zNZHO_Sjf+solana-sol"
$ws.Range("D12").Value = "'22.52"
$ws.Range("E12").Value = "  +2.84%  "

# Row 13
$ws.Range("B13").Value = "WrappedEther"
$ws.Range("C13").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D13").Value = "1.973.90"
$ws.Range("E13").Value = "  +7.26%  "

# Row 14
$ws.Range("B14").Value = "Chainlink"
$ws.Range("C14").Value = "https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link"
$ws.Range("D14").Value = "'7.214"
$ws.Range("E14").Value = "  +2.51%  "

# Row 15
$ws.Range("B15").Value = "Polkadot"
$ws.Range("C15").Value = "https://coinranking.com/coin/25W7FG7om+polkadot-dot"
$ws.Range("D15").Value = "'5.852"
$ws.Range("E15").Value = "  +2.74%  "

# Row 16
$ws.Range("B16").Value = "TRON"
$ws.Range("C16").Value = "https://coinranking.com/coin/qUhEFk1I61atv+tron-trx"
$ws.Range("D16").Value = "'0.07103"
$ws.Range("E16").Value = "  +2.33%  "

# Row 17
$ws.Range("B17").Value = "Litecoin"
$ws.Range("C17").Value = "https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc"
$ws.Range("D17").Value = "'88.80"
$ws.Range("E17").Value = "  +0.81%  "

# Row 18
$ws.Range("B18").Value = "BinanceUSD"
$ws.Range("C18").Value = "https://coinranking.com/coin/vSo2fu9iE1s0Y+binanceusd-busd"
$ws.Range("D18").Value = "'1.000"
$ws.Range("E18").Value = "  -0.29%  "

# Row 19
$ws.Range("B19").Value = "ShibaInu"
$ws.Range("C19").Value = "https://coinranking.com/coin/xz24e0BjL+shibainu-shib"
$ws.Range("D19").Value = "'0.000009973"
$ws.Range("E19").Value = "  -0.21%  "

# Row 20
$ws.Range("B20").Value = "Avalanche"
$ws.Range("C20").Value = "https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax"
$ws.Range("D20").Value = "'17.20"
$ws.Range("E20").Value = "  +1.07%  "

# Row 21
$ws.Range("B21").Value = "Dai"
$ws.Range("C21").Value = "https://coinranking.com/coin/MoTuySvg7+dai-dai"
$ws.Range("D21").Value = "'0.9977"
$ws.Range("E21").Value = "  -0.23%  "

# Row 22
$ws.Range("B22").Value = "WrappedBTC"
$ws.Range("C22").Value = "https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc"
$ws.Range("D22").Value = "29.550.43"
$ws.Range("E22").Value = "  +2.14%  "

# Row 23
$ws.Range("B23").Value = "Uniswap"
$ws.Range("C23").Value = "https://coinranking.com/coin/_H5FVG9iW+uniswap-uni"
$ws.Range("D23").Value = "'5.575"
$ws.Range("E23").Value = "  +5.26%  "

# Row 24
$ws.Range("B24").Value = "Cosmos"
$ws.Range("C24").Value = "https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom"
$ws.Range("D24").Value = "'11.25"
$ws.Range("E24").Value = "  +2.47%  "

# Row 25
$ws.Range("B25").Value = "WrappedliquidstakedEther2.0"
$ws.Range("C25").Value = "https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth"
$ws.Range("D25").Value = "2.217.97"
$ws.Range("E25").Value = "  +6.73%  "

# Row 26
$ws.Range("B26").Value = "Toncoin"
$ws.Range("C26").Value = "https://coinranking.com/coin/67YlI0K1b+toncoin-ton"
$ws.Range("D26").Value = "'2.119"
$ws.Range("E26").Value = "  +3.22%  "

# Row 27
$ws.Range("B27").Value = "Monero"
$ws.Range("C27").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D27").Value = "'157.85"
$ws.Range("E27").Value = "  +1.16%  "

# Row 28
$ws.Range("B28").Value = "EthereumClassic"
$ws.Range("C28").Value = "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
$ws.Range("D28").Value = "'19.65"
$ws.Range("E28").Value = "  +1.80%  "

# Row 29
$ws.Range("B29").Value = "InternetComputer(DFINITY)"
$ws.Range("C29").Value = "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
$ws.Range("D29").Value = "'6.014"
$ws.Range("E29").Value = "  +1.18%  "

# Row 30
$ws.Range("B30").Value = "BitcoinCash"
$ws.Range("C30").Value = "https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch"
$ws.Range("D30").Value = "'120.49"
$ws.Range("E30").Value = "  +2.35%  "

# Row 31
$ws.Range("B31").Value = "LidoDAOToken"
$ws.Range("C31").Value = "https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo"
$ws.Range("D31").Value = "'1.919"
$ws.Range("E31").Value = "  -0.55%  "

# Row 32
$ws.Range("B32").Value = "Stellar"
$ws.Range("C32").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$ws.Range("D32").Value = "'0.09410"
$ws.Range("E32").Value = "  +0.83%  "

# Row 33
$ws.Range("B33").Value = "ImmutableX"
$ws.Range("C33").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D33").Value = "'0.8947"
$ws.Range("E33").Value = "  -1.73%  "

# Row 34
$ws.Range("B34").Value = "Filecoin"
$ws.Range("C34").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("D34").Value = "'5.278"
$ws.Range("E34").Value = "  -0.36%  "

# Row 35
$ws.Range("B35").Value = "ARBITRUM"
$ws.Range("C35").Value = "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
$ws.Range("D35").Value = "'1.348"
$ws.Range("E35").Value = "  +1.45%  "

# Row 36
$ws.Range("B36").Value = "PEPE"
$ws.Range("C36").Value = "https://coinranking.com/coin/03WI8NQPF+pepe-pepe"
$ws.Range("D36").Value = "'0.000004189"
$ws.Range("E36").Value = "  +155.69%  "

# Row 37
$ws.Range("B37").Value = "HuobiToken"
$ws.Range("C37").Value = "https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht"
$ws.Range("D37").Value = "'3.167"
$ws.Range("E37").Value = "  -2.99%  "

# Row 38
$ws.Range("B38").Value = "Hedera"
$ws.Range("C38").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D38").Value = "'0.05828"
$ws.Range("E38").Value = "  +0.90%  "

# Row 39
$ws.Range("B39").Value = "TrustWalletToken"
$ws.Range("C39").Value = "https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt"
$ws.Range("D39").Value = "'1.177"
$ws.Range("E39").Value = "  -2.05%  "

# Row 40
$ws.Range("B40").Value = "VeChain"
$ws.Range("C40").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D40").Value = "'0.02129"
$ws.Range("E40").Value = "  +3.03%  "

# Row 41
$ws.Range("B41").Value = "FraxShare"
$ws.Range("C41").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$ws.Range("D41").Value = "'7.937"
$ws.Range("E41").Value = "  +3.84%  "

# Row 42
$ws.Range("D42").Value = "'0.9978"
$ws.Range("E42").Value = "  -0.36%  "

# Row 43
$ws.Range("B43").Value = "TheSandbox"
$ws.Range("C43").Value = "https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand"
$ws.Range("D43").Value = "'0.5775"
$ws.Range("E43").Value = "  +1.76%  "

# Row 44
$ws.Range("B44").Value = "Algorand"
$ws.Range("C44").Value = "https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo"
$ws.Range("D44").Value = "'0.1820"
$ws.Range("E44").Value = "  +2.67%  "

# Row 45
$ws.Range("B45").Value = "Aptos"
$ws.Range("C45").Value = "https://coinranking.com/coin/HGYj5JCv5+aptos-apt"
$ws.Range("D45").Value = "'9.813"
$ws.Range("E45").Value = "  +0.61%  "

# Row 46
$ws.Range("B46").Value = "EnergySwap"
$ws.Range("C46").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D46").Value = "'12.08"
$ws.Range("E46").Value = "  +0.77%  "

# Row 47
$ws.Range("B47").Value = "Decentraland"
$ws.Range("C47").Value = "https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana"
$ws.Range("D47").Value = "'0.5372"
$ws.Range("E47").Value = "  +0.11%  "

# Row 48
$ws.Range("B48").Value = "RenderToken"
$ws.Range("C48").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D48").Value = "'2.187"
$ws.Range("E48").Value = "  -5.16%  "

# Row 49
$ws.Range("B49").Value = "MXToken"
$ws.Range("C49").Value = "https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx"
$ws.Range("D49").Value = "'2.642"
$ws.Range("E49").Value = "  +4.88%  "

# Row 50
$ws.Range("B50").Value = "Cronos"
$ws.Range("C50").Value = "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
$ws.Range("D50").Value = "'0.06981"
$ws.Range("E50").Value = "  -0.92%  "

# Row 51
$ws.Range("B51").Value = "NEARProtocol"
$ws.Range("C51").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("D51").Value = "'1.871"
$ws.Range("E51").Value = "  +0.93%  "
